# Equipment_Record.xlsx - "Added X310 USRP Kit"
#
# Context: row 2 column A previously said the generic "USRP", but the sheet
# is adding a second USRP kit entry (X310) further down, so the original
# item is re-labelled "USRP B210" to distinguish it, and a new row 24 is
# appended describing the USRP X310 kit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new USRP X310 kit row (Item cell first)...
$ws.Range("A24").Value = "USRP X310"

# ...then disambiguate the existing B210 kit's "Item" cell...
$ws.Range("A2").Value = "USRP B210"

# ...then finish filling in the rest of the new row.
$ws.Range("B24").Value = "USRP X310 KIT (KINTEX7-410T FPGA, 2 CHANNELS, 10GIGE AND PCIE BUS)"
$ws.Range("C24").Value = "PA1289856"
$ws.Range("E24").Value = "San Diego - Palomar"

# Match the plain-black font styling (no theme color) used by the rest of
# the "Item" / "Physical Location" columns in this table.
$ws.Range("A24").Font.Color = 0
$ws.Range("E24").Font.Color = 0

# Leave the new row selected, same as the author's last action.
$ws.Range("E25").Select()
